# Riaz.docx -> two new paragraphs about "Galib" replace the original
# "Riaz Khan" bio sentences.
#
# The edit:
#   1. Merges the document's two paragraphs into a single paragraph.
#   2. Replaces the text with "Galib is a businessman. He work in the
#      town of Badarganj. "
#   3. Wraps the word "Badarganj" with <w:proofErr w:type="spellStart"/>
#      / <w:proofErr w:type="spellEnd"/> (the markers Word itself adds
#      around a word its spell-checker doesn't recognize), which means
#      that word ends up in its own run, flanked by the "He work in the
#      town of " run and the ". " run.

$d = $word.ActiveDocument

# Locate (and sanity-check) the original two sentences across the
# paragraph mark that separates them ("^p" = Word's paragraph-mark find
# code) before touching anything.
$probe = $d.Content
$hasOriginal = $probe.Find.Execute(
    "Riaz Khan is a Software Engineer. Riaz Khan work in MyMedicalHUB International Ltd for 5 years.^p" + `
    "Riaz Khan completed graduation from RUET.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $hasOriginal) {
    throw "Could not locate the expected original paragraphs in the document."
}

# Build the replacement body: a single paragraph, split into three runs
# so "Badarganj" can be flagged with proofErr spell-check markers, just
# like Word does automatically for a word outside its dictionary.
$newBodyXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="66D8C7DE" w14:textId="074B2852" w:rsidR="00A4269E" w:rsidRDefault="00F256C9"><w:r><w:t xml:space="preserve">Galib is a businessman. He work in the town of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Badarganj</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# Replacing the whole-document Content range with this fragment both
# substitutes the text AND merges the two original paragraphs into one
# (the fragment only contains a single <w:p>).
$d.Content.InsertXML($newBodyXml)
